# Updated symbol list on Sat Dec 24 22:36:01 UTC 2022 with GitHub Actions
# Refresh the Price (column D) and a couple of Volume(1h) (column E) values
# for the cryptos sheet. The Price column stores numeric-looking values as
# text, so force NumberFormat = "@" (Text) on each touched cell before
# assigning the new value - otherwise Excel's COM type-inference would
# silently convert the string into a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2"  "244.52"
Set-TextValue "D3"  "21.89"
Set-TextValue "D4"  "5.393"
Set-TextValue "D5"  "0.06037"
Set-TextValue "D7"  "0.8143"
Set-TextValue "D8"  "0.9241"
Set-TextValue "D9"  "0.1441"
Set-TextValue "D10" "0.07460"
Set-TextValue "D11" "0.03391"
Set-TextValue "D12" "0.03050"
Set-TextValue "D13" "0.09425"
Set-TextValue "D14" "4.006"
Set-TextValue "D15" "0.001597"
Set-TextValue "D16" "0.04805"
Set-TextValue "D17" "0.0005945"
Set-TextValue "D18" "0.005510"
Set-TextValue "D19" "0.004153"
Set-TextValue "D20" "0.0009898"
Set-TextValue "D22" "6.432"
Set-TextValue "D26" "0.00008508"
Set-TextValue "D27" "0.0002903"
Set-TextValue "D40" "0.04003"
Set-TextValue "D41" "0.006416"

$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

Set-TextValue "D42" "0.1073"
Set-TextValue "D43" "0.002903"
Set-TextValue "D45" "0.00005252"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Write-Host "Symbol list updated"
